$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..13 (Title, Company, Location, Match Score, Matched Keywords, Posted At, Apply Link)
$data = @(
    @("AI Data Analyst", "Maximus", "Tysons, VA, US USA", 20, "Generative AI, RAG, Copilot, TensorFlow, PyTorch, XGBoost, Keras, AWS SageMaker, Snowflake, Databricks", "2026-02-24", "https://www.indeed.com/viewjob?jk=5908ca8b0f7aacbd"),
    @("Data Engineering Intern", "Global Partners LP", "Waltham, MA, US USA", 16.7, "Data Scientist, RAG, S3, Glue, Redshift, BigQuery, Docker, CI/CD, Git, Snowflake", "2026-02-24", "https://www.indeed.com/viewjob?jk=b60d9899ce76c629"),
    @("Data Engineer", "Books-A-Million", "Birmingham, AL, US USA", 14.4, "RAG, S3, Redshift, BigQuery, Data Lake, Git, BigQuery, Redshift, Python, SQL", "2026-02-23", "https://www.indeed.com/viewjob?jk=7bb1becaa3bb4909"),
    @("HPC Software Engineer", "Cadre5", "Knoxville, TN, US USA", 13.3, "RAG, TensorFlow, PyTorch, Keras, Docker, Kubernetes, Git, Dask, Python, R", "2026-02-23", "https://www.indeed.com/viewjob?jk=51ad2b47ec2abe67"),
    @("Software Developer 4 - Redwood City, CA (in-office role)", "Oracle", "Redwood City, CA, US USA", 12.2, "RAG, XGBoost, MLflow, FastAPI, Docker, Kubernetes, CI/CD, PySpark, R, Scala", "2026-02-23", "https://www.indeed.com/viewjob?jk=701b0e63a9a8e050"),
    @("Cloud Engineer II", "GSK", "Philadelphia, PA, US USA", 12.2, "RAG, Docker, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Python, R, Java", "2026-02-19", "https://www.indeed.com/viewjob?jk=3258391b26d34a0f"),
    @("Senior Cloud Engineer", "GSK", "Philadelphia, PA, US USA", 11.1, "RAG, Docker, Kubernetes, Jenkins, GitHub Actions, Terraform, Git, Python, R, Scala", "2025-12-18", "https://www.indeed.com/viewjob?jk=7ed37c67470dc51b"),
    @("Sr. Site Reliability Engineer - SRE", "AT&T", "Dallas, TX, US USA", 10, "RAG, Docker, Kubernetes, Kafka, Python, SQL, R, Java, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=bd62853e90f5d50a"),
    @("Sr. Site Reliability Engineer - SRE", "AT&T", "Alpharetta, GA, US USA", 10, "RAG, Docker, Kubernetes, Kafka, Python, SQL, R, Java, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=326d5680ae095aba"),
    @("Analytics Engineer II", "CDW", "Remote, US USA", 10, "Data Scientist, Data Lake, AKS, Git, Power BI, Python, SQL, R, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=3133b5a21a1473ac"),
    @("Entry-Level AI / Machine Learning Software Engineer", "MTSI", "Huntsville, AL, US USA", 10, "RAG, TensorFlow, PyTorch, Docker, Git, Matplotlib, Seaborn, Python, R", "2026-02-23", "https://www.indeed.com/viewjob?jk=6cd62aed0ccfaf4c"),
    @("Associate Data Engineer", "CVS Health", "TX, US USA", 10, "Data Scientist, RAG, BigQuery, Dataflow, Apache Airflow, BigQuery, Python, SQL, R", "2026-01-06", "https://www.indeed.com/viewjob?jk=9291f3f30c723313")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    # Prefix with an apostrophe so Excel keeps this date-looking string as plain text
    # (matches the source data, which stores dates as text, not date serials).
    $ws.Cells.Item($r, 6).Value = "'" + $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Remove the now-obsolete last row (row 14), shifting the dimension back to A1:G13
$ws.Rows.Item(14).Delete()
